# Weekly update for "Hortaliza, Vega Monumental Concepción - Acelga":
# a new week's pair of price rows (Primera / Segunda) is inserted at the
# top of the data block (row 96), pushing the existing historical rows
# down by two and extending the sheet from A1:R151 to A1:R153.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the current row 96 - this shifts rows
# 96:151 down to 98:153 and keeps per-row formatting (e.g. the date style
# on column D) intact.
$ws.Range("A96:R97").EntireRow.Insert()

# Row 96: "Primera" quality entry for the new week (serial 44460 = 2021-09-21)
$ws.Cells.Item(96, 1).Value = 11
$ws.Cells.Item(96, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(96, 3).Value = "Bíobío"
$ws.Cells.Item(96, 4).Value = 44460
$ws.Cells.Item(96, 5).Value = 8
$ws.Cells.Item(96, 6).Value = 100112009
$ws.Cells.Item(96, 7).Value = "Acelga"
$ws.Cells.Item(96, 8).Value = "Sin especificar"
$ws.Cells.Item(96, 9).Value = "Primera"
$ws.Cells.Item(96, 10).Value = 200
$ws.Cells.Item(96, 11).Value = 600
$ws.Cells.Item(96, 12).Value = 700
$ws.Cells.Item(96, 13).Value = 650
$ws.Cells.Item(96, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(96, 15).Value = "Región de Ñuble"
$ws.Cells.Item(96, 16).Value = 650
$ws.Cells.Item(96, 17).Value = 1
$ws.Cells.Item(96, 18).Value = "Hortaliza"

# Row 97: "Segunda" quality entry for the new week
$ws.Cells.Item(97, 1).Value = 11
$ws.Cells.Item(97, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(97, 3).Value = "Bíobío"
$ws.Cells.Item(97, 4).Value = 44460
$ws.Cells.Item(97, 5).Value = 8
$ws.Cells.Item(97, 6).Value = 100112009
$ws.Cells.Item(97, 7).Value = "Acelga"
$ws.Cells.Item(97, 8).Value = "Sin especificar"
$ws.Cells.Item(97, 9).Value = "Segunda"
$ws.Cells.Item(97, 10).Value = 100
$ws.Cells.Item(97, 11).Value = 500
$ws.Cells.Item(97, 12).Value = 500
$ws.Cells.Item(97, 13).Value = 500
$ws.Cells.Item(97, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(97, 15).Value = "Región de Ñuble"
$ws.Cells.Item(97, 16).Value = 500
$ws.Cells.Item(97, 17).Value = 1
$ws.Cells.Item(97, 18).Value = "Hortaliza"
